$d = $word.ActiveDocument

$replacements = @(
    @("15×27=", "56×94="),
    @("13×40=", "32×62="),
    @("72×25=", "46×46="),
    @("31×53=", "34×30="),
    @("59×81=", "25×44="),
    @("78×66=", "31×32="),
    @("15×94=", "62×59="),
    @("81×89=", "12×33="),
    @("54×32=", "72×38="),
    @("21×95=", "51×69="),
    @("95×68=", "41×74="),
    @("87×17=", "89×15="),
    @("34×80=", "22×64="),
    @("98×16=", "56×89="),
    @("16×52=", "60×56="),
    @("54×31=", "47×57="),
    @("24×32=", "29×98="),
    @("26×29=", "61×35="),
    @("36×74=", "24×15="),
    @("54×56=", "42×58="),
    @("48×17=", "49×25="),
    @("63×65=", "87×30="),
    @("28×12=", "62×79="),
    @("59×77=", "28×44="),
    @("99×47=", "34×94=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
